$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Readout values for the two "Weight" log rows so they
# reflect the new naming scheme used in this test run.
$ws.Range("E4").Value = " Weight MeOH_2"
$ws.Range("E2").Value = " Weight 2MIM_2"

# Widen column E to fit the new (longer) text and match the manual
# "AutoFit" that was performed after editing the values (stored column
# width of 22 once Excel applies its internal padding).
$ws.Columns.Item(5).ColumnWidth = 21.14

# Move the active selection to where the user left off after editing.
$ws.Range("F9").Select()
